# Update TPM-derived NATMI metrics for the Pdgfb-Lrp1 ligand-receptor pair.
# Ligand-expressing-cell count (E) moved from 2 to 3 for ECs/FAPs source rows,
# which cascades through the ligand/receptor/edge specificity + weight columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 46.11811066666667
$ws.Range("H2").Value = 138.354332
$ws.Range("I2").Value = 0.95896098489411
$ws.Range("J2").Value = 0.9589609848941099
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 159.3964271360302
$ws.Range("R2").Value = 1434.567844224272
$ws.Range("S2").Value = 0.009437648871001066
$ws.Range("T2").Value = 0.009437648871001065

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 46.11811066666667
$ws.Range("H3").Value = 138.354332
$ws.Range("I3").Value = 0.95896098489411
$ws.Range("J3").Value = 0.9589609848941099
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 13909.22903329245
$ws.Range("R3").Value = 125183.0612996321
$ws.Range("S3").Value = 0.8235468136969
$ws.Range("T3").Value = 0.8235468136968997

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 46.11811066666667
$ws.Range("H4").Value = 138.354332
$ws.Range("I4").Value = 0.95896098489411
$ws.Range("J4").Value = 0.9589609848941099
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 2127.67054976163
$ws.Range("R4").Value = 19149.03494785467
$ws.Range("S4").Value = 0.1259765223262091
$ws.Range("T4").Value = 0.125976522326209

# Row 5
$ws.Range("I5").Value = 0.002799731840346333
$ws.Range("J5").Value = 0.002799731840346333
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 0.4653653895413333
$ws.Range("R5").Value = 4.188288505871999
$ws.Range("S5").Value = 0.00002755366115866327
$ws.Range("T5").Value = 0.00002755366115866326

# Row 6
$ws.Range("I6").Value = 0.002799731840346333
$ws.Range("J6").Value = 0.002799731840346333
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("Q6").Value = 40.60865041707466
$ws.Range("R6").Value = 365.4778537536719
$ws.Range("S6").Value = 0.002404383778530448
$ws.Range("T6").Value = 0.002404383778530448

# Row 7
$ws.Range("I7").Value = 0.002799731840346333
$ws.Range("J7").Value = 0.002799731840346333
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("Q7").Value = 6.211834556118667
$ws.Range("S7").Value = 0.0003677944006572218
$ws.Range("T7").Value = 0.0003677944006572217

# Row 8
$ws.Range("I8").Value = 0.0382392832655437
$ws.Range("J8").Value = 0.0382392832655437
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 6.356051210408
$ws.Range("R8").Value = 57.20446089367201
$ws.Range("S8").Value = 0.0003763332755177714
$ws.Range("T8").Value = 0.0003763332755177713

# Row 9
$ws.Range("I9").Value = 0.0382392832655437
$ws.Range("J9").Value = 0.0382392832655437
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.03283954236664705
$ws.Range("T9").Value = 0.03283954236664704

# Row 10
$ws.Range("I10").Value = 0.0382392832655437
$ws.Range("J10").Value = 0.0382392832655437
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.005023407623378888
$ws.Range("T10").Value = 0.005023407623378887
